# Update stats for 2025-12 (row 25 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B25").Value = 6443
$ws.Range("D25").Value = 6004020
$ws.Range("E25").Value = 931.8671426354183
$ws.Range("F25").Value = 9.370225768120854
$ws.Range("H25").Value = 25.74269410992522
